$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "MEC-3B-Comam. Pneumáticos"
$ws.Range("D2").Value = "MEC-3B-Coman. Hidráulicos"
$ws.Range("E2").Value = "MCT-3A-Eletrohidráulica"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "MEC-3B-Comam. Pneumáticos"
$ws.Range("D3").Value = "MEC-3B-Comam. Pneumáticos"
$ws.Range("E3").Value = "MCT-3A-Eletrohidráulica"

# Row 4
$ws.Range("D4").Value = "MEC-3B-Coman. Hidráulicos"
$ws.Range("E4").Value = "-"

# Row 6
$ws.Range("D6").Value = "MEC-3B-Coman. Hidráulicos"
$ws.Range("E6").Value = "-"

# Row 7
$ws.Range("D7").Value = "MEC-3B-Coman. Hidráulicos"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "MCT-3A-Eletrohidráulica"

# Row 8
$ws.Range("D8").Value = "MEC-3B-Comam. Pneumáticos"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "MCT-3A-Eletrohidráulica"

# Row 10
$ws.Range("B10").Value = "MEC-3A-Cont. Lóg. Prog. CLP"
$ws.Range("C10").Value = "MEC-2A-Maq. Term. FLuxo"
$ws.Range("E10").Value = "MEC-3A-Cont. Lóg. Prog. CLP"

# Row 11
$ws.Range("C11").Value = "MEC-2A-Maq. Term. FLuxo"
$ws.Range("E11").Value = "MEC-3A-Camam. Hidráulicos"

# Row 12
$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "-"

# Row 14
$ws.Range("C14").Value = "MEC-3A-Camam. Hidráulicos"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"

# Row 15
$ws.Range("C15").Value = "MEC-3A-Cont. Lóg. Prog. CLP"
$ws.Range("E15").Value = "-"

# Row 16
$ws.Range("B16").Value = "MEC-3A-Camam. Hidráulicos"
$ws.Range("C16").Value = "MEC-3A-Camam. Hidráulicos"
$ws.Range("E16").Value = "MEC-3A-Cont. Lóg. Prog. CLP"
